$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DAC80508")

# Shared-string allocation order observed in the target file: WB_IN_1, WB_IN_0,
# WB_OUT_0, WB_OUT_1 -- so write column A for row 19 before row 18.
$ws.Range("A19").Value = "WB_IN_1"
$ws.Range("A18").Value = "WB_IN_0"
$ws.Range("A20").Value = "WB_OUT_0"
$ws.Range("A21").Value = "WB_OUT_1"

$ws.Range("B18").Value = "0x03"
$ws.Range("C18").Value = "0x00000000"
$ws.Range("D18").Value = 32
$ws.Range("E18").Value = "None"
$ws.Range("F18").Value = "None"

$ws.Range("B19").Value = "0x04"
$ws.Range("C19").Value = "0x00000000"
$ws.Range("D19").Value = 32
$ws.Range("E19").Value = "None"
$ws.Range("F19").Value = "None"

$ws.Range("B20").Value = "0x22"
$ws.Range("C20").Value = "0x00000000"
$ws.Range("D20").Value = 32
$ws.Range("E20").Value = "None"
$ws.Range("F20").Value = "None"

$ws.Range("B21").Value = "0x23"
$ws.Range("C21").Value = "0x00000000"
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = "None"
$ws.Range("F21").Value = "None"

$ws.Activate() | Out-Null
$ws.Range("A22").Select() | Out-Null
